$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Results")
$tbl = $ws.ListObjects.Item("Table1")

# Add a new table column ("Test type") - this expands the table/autofilter range to column I
$newCol = $tbl.ListColumns.Add()

# Fill in the header and the per-row test-type descriptions, in the same order
# the values were first authored (controls shared-string insertion order).
$ws.Range("I25").Value = "Semaphore creation and deletion API"
$ws.Range("I26").Value = "atomSemGet() and atomSemPut"
$ws.Range("I27").Value = "Semaphore priority and FIFO on waking threads"
$ws.Range("I28").Value = "Synchronisation between threads"
$ws.Range("I29").Value = "Semaphore stress test Get and Put"
$ws.Range("I32").Value = "Semaphore stress test Get and Put"
$ws.Range("I30").Value = "Semaphore basic counting test"
$ws.Range("I31").Value = "Semaphore for basic mutual exclusion test"
$ws.Range("I33").Value = "Semaphore deletion API with multiple blocked threads"
$ws.Range("I34").Value = "atomTimerDelay() API test"
$ws.Range("I35").Value = "atomTimerDelay() with 3 threads test"
$ws.Range("I36").Value = "atomTimerRegister() API test"
$ws.Range("I37").Value = "atomTimerRegister() API test"
$ws.Range("I38").Value = "atomTimerCancel() API test"
$ws.Range("I39").Value = "atomTimerCancel() API test"
$ws.Range("I40").Value = "timer subsystem behaviour test"
$ws.Range("I41").Value = "timer register within a timer callback test"
$ws.Range("I3").Value = "Basic context-switch test"
$ws.Range("I2").Value = "Bad parameter handling of public API test"
$ws.Range("I4").Value = "Scheduling with different priorities and preemption test"
$ws.Range("I5").Value = "Round-Robin timeslicing test"
$ws.Range("I14").Value = "Timeouts on mutex test"
$ws.Range("I6").Value = "Mutex creation and deletion API test"
$ws.Range("I7").Value = "atomMutexGet() and atomMutexPUT() API test"
$ws.Range("I8").Value = "Stress test the mutex Get and Put operations"
$ws.Range("I9").Value = "Mutex priority and FIFO on waking threads test"
$ws.Range("I10").Value = "Basic mutex usage test"
$ws.Range("I11").Value = "Mutex lock count test"
$ws.Range("I12").Value = "Mutex ownership test"
$ws.Range("I13").Value = "Mutex deletion API test"
$ws.Range("I15").Value = "Queue creation and deletion API test"
$ws.Range("I17").Value = "Queue deletion while threads blocking on atomQueueGet() test"
$ws.Range("I18").Value = "Queue deletion while threads blocking on atomQueuePut() test"
$ws.Range("I19").Value = "atomQueueGet() and atomQueuePut() stress test"
$ws.Range("I20").Value = "Queue priority and FIFO on waking threads test"
$ws.Range("I16").Value = "Queue basic operation test"
$ws.Range("I21").Value = "Queue basic operation test"
$ws.Range("I22").Value = "Queue deletion API test"
$ws.Range("I23").Value = "Timeout on queue test"
$ws.Range("I24").Value = "Queue get and put stress test"
$ws.Range("I1").Value = "Test type"

# Left-align the data cells in the new column (header stays default-aligned).
$ws.Range("I2:I41").HorizontalAlignment = -4131

# Widen the new column to fit the long descriptions.
$ws.Columns.Item(9).ColumnWidth = 52.75
